# Precolombinos.docx — "se incluye contraseña y recuperacion de contraseña"
#
# The document is a bulleted file/folder tree (list style "Prrafodelista",
# numId=1). This change adds the new "producto" (product) and "venta"
# (sale) files that were introduced alongside the password / password-
# recovery feature, at the same list level (ilvl=2) as their siblings:
#
#   controllers: + productoController.js, ventaController.js   (after authController.js)
#   models:      + productoModel.js                            (after userModel.js)
#   routes:      + productoRoutes.js, ventaRoutes.js,
#                  usuarioRoutes.js                             (after authRoutes.js)
#   js:          + productos.js (before register.js),
#                  venta.js (after register.js)
#   pages:       + productos.html (before register.html),
#                  reportes.html, ventas.html (after register.html)

$d = $word.ActiveDocument

function Get-ParaByText($doc, $text) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $p
        }
    }
    throw "Paragraph not found: $text"
}

# Insert a brand-new list paragraph right after the paragraph whose text is
# $anchorText, carrying $newText. The new paragraph inherits the anchor's
# pPr (style / numPr / ilvl), matching what Word itself does when you press
# Enter at the end of a list line.
function Insert-ListParaAfter($doc, $anchorText, $newText) {
    $anchor = Get-ParaByText $doc $anchorText
    $anchor.Range.InsertParagraphAfter() | Out-Null
    $created = Get-ParaByText $doc $anchorText
    $created.Next().Range.Text = $newText
}

# Insert a brand-new list paragraph right before the paragraph whose text is
# $anchorText, carrying $newText, inheriting the anchor's pPr.
function Insert-ListParaBefore($doc, $anchorText, $newText) {
    $anchor = Get-ParaByText $doc $anchorText
    $anchor.Range.InsertParagraphBefore() | Out-Null
    $created = Get-ParaByText $doc $anchorText
    $created.Previous().Range.Text = $newText
}

# --- backend/controllers -------------------------------------------------
Insert-ListParaAfter $d "authController.js" "productoController.js"
Insert-ListParaAfter $d "productoController.js" "ventaController.js"

# --- backend/models --------------------------------------------------------
Insert-ListParaAfter $d "userModel.js" "productoModel.js"

# --- backend/routes ----------------------------------------------------
Insert-ListParaAfter $d "authRoutes.js" "productoRoutes.js"
Insert-ListParaAfter $d "productoRoutes.js" "ventaRoutes.js"
Insert-ListParaAfter $d "ventaRoutes.js" "usuarioRoutes.js"

# --- frontend/js ---------------------------------------------------------
Insert-ListParaBefore $d "register.js" "productos.js"
Insert-ListParaAfter $d "register.js" "venta.js"

# --- frontend/pages --------------------------------------------------------
Insert-ListParaBefore $d "register.html" "productos.html"
Insert-ListParaAfter $d "register.html" "reportes.html"
Insert-ListParaAfter $d "reportes.html" "ventas.html"

Write-Host "Paragraphs after edit:" $d.Paragraphs.Count
